$d = $word.ActiveDocument

# Each (old, new) pair below is applied via Find/Replace against the whole
# document content. The pairs are ordered exactly as they occur in the
# document so that no newly-inserted text is accidentally re-matched by a
# later replacement (this matters because "97÷2=" is both a source value
# and the replacement target produced from "35÷7=").
$pairs = @(
    @("56÷5=", "36÷2="),
    @("56÷4=", "63÷2="),
    @("92÷4=", "89÷5="),
    @("30÷3=", "20÷8="),
    @("10÷3=", "52÷6="),
    @("62÷4=", "25÷8="),
    @("14÷9=", "44÷2="),
    @("52÷2=", "86÷4="),
    @("22÷5=", "90÷6="),
    @("97÷2=", "95÷7="),
    @("90÷3=", "82÷6="),
    @("25÷5=", "33÷2="),
    @("65÷2=", "45÷5="),
    @("51÷9=", "86÷6="),
    @("51÷6=", "86÷8="),
    @("38÷9=", "20÷2="),
    @("54÷2=", "26÷6="),
    @("73÷4=", "52÷5="),
    @("90÷5=", "52÷6="),
    @("69÷3=", "95÷2="),
    @("95÷3=", "15÷8="),
    @("37÷9=", "34÷5="),
    @("35÷7=", "97÷2="),
    @("80÷7=", "29÷8="),
    @("60÷9=", "36÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
